$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.200.21'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5191'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06268'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07789'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.472'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.657.23'
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('D14').Value = '1.886.80'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5463'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '0.0₅8117'
$ws.Range('E16').Value = '  -1.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '26.208.30'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.607'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('E22').Value = '  -2.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.45%  '
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.07'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('E26').Value = '  -3.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.284'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.438'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05944'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.274'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.547'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.266'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.583'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9595'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.419'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.769'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5690'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.78%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01592'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.988'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8487'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.44'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '1.001.85'
$ws.Range('E44').Value = '  -8.65%  '
$ws.Range('D45').Value = '1.801.28'
$ws.Range('E46').Value = '  +8.38%  '
$ws.Range('E47').Value = '  -2.69%  '
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.020'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4336'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05161'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.82%  '
